$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 0.9168783282298989
$ws.Cells.Item(2, 3).Value = 0.1103637748559407
$ws.Cells.Item(2, 4).Value = 0.121236058883909
$ws.Cells.Item(2, 5).Value = 0.1322696356981972
$ws.Cells.Item(2, 6).Value = 1.940048721761194
$ws.Cells.Item(2, 10).Value = 0.1762184897146213
$ws.Cells.Item(2, 11).Value = 0.4782845967818332
$ws.Cells.Item(2, 12).Value = 0.263171286780711
$ws.Cells.Item(2, 14).Value = 2.252590438682947
$ws.Cells.Item(2, 15).Value = 5.132683292418307

$ws.Cells.Item(3, 2).Value = 0.8797396025435376
$ws.Cells.Item(3, 3).Value = 0.1093442579316104
$ws.Cells.Item(3, 4).Value = 0.1194899340974089
$ws.Cells.Item(3, 5).Value = 0.1323217267412211
$ws.Cells.Item(3, 6).Value = 1.946675811283775
$ws.Cells.Item(3, 10).Value = 0.1771558591426654
$ws.Cells.Item(3, 11).Value = 0.4443327767251333
$ws.Cells.Item(3, 12).Value = 0.2593479686275089
$ws.Cells.Item(3, 14).Value = 2.27298144758604
$ws.Cells.Item(3, 15).Value = 5.157466130098129

$ws.Cells.Item(4, 2).Value = 0.8572861856528959
$ws.Cells.Item(4, 3).Value = 0.1087140804762328
$ws.Cells.Item(4, 4).Value = 0.1184611853798785
$ws.Cells.Item(4, 5).Value = 0.1323975060883154
$ws.Cells.Item(4, 6).Value = 1.951551618945111
$ws.Cells.Item(4, 10).Value = 0.1777840706764238
$ws.Cells.Item(4, 11).Value = 0.4236228971695937
$ws.Cells.Item(4, 12).Value = 0.2571021721292581
$ws.Cells.Item(4, 14).Value = 2.286138749727554
$ws.Cells.Item(4, 15).Value = 5.174836045806416

$ws.Cells.Item(5, 2).Value = 0.8482249817101319
$ws.Cells.Item(5, 3).Value = 0.1084562410388159
$ws.Cells.Item(5, 4).Value = 0.1180529336918141
$ws.Cells.Item(5, 5).Value = 0.1324394302075849
$ws.Cells.Item(5, 6).Value = 1.953741668300609
$ws.Cells.Item(5, 10).Value = 0.1780533363990404
$ws.Cells.Item(5, 11).Value = 0.4152183855656375
$ws.Cells.Item(5, 12).Value = 0.2562126944634073
$ws.Cells.Item(5, 14).Value = 2.291660731319875
$ws.Cells.Item(5, 15).Value = 5.182456355751597

$ws.Cells.Item(6, 2).Value = 0.8467257580255705
$ws.Cells.Item(6, 3).Value = 0.1084133648533125
$ws.Cells.Item(6, 4).Value = 0.117985808354149
$ws.Cells.Item(6, 5).Value = 0.1324470595755045
$ws.Cells.Item(6, 6).Value = 1.954117599665928
$ws.Cells.Item(6, 10).Value = 0.1780988495499329
$ws.Cells.Item(6, 11).Value = 0.4138249492640398
$ws.Cells.Item(6, 12).Value = 0.2560665532878161
$ws.Cells.Item(6, 14).Value = 2.292587334124468
$ws.Cells.Item(6, 15).Value = 5.183754448296355

$ws.Cells.Item(7, 2).Value = 0.8571636227548254
$ws.Cells.Item(7, 3).Value = 0.1087106073414432
$ws.Cells.Item(7, 4).Value = 0.1184556350450947
$ws.Cells.Item(7, 5).Value = 0.132398026736297
$ws.Cells.Item(7, 6).Value = 1.951580331961082
$ws.Cells.Item(7, 10).Value = 0.1777876483515044
$ws.Cells.Item(7, 11).Value = 0.4235094087985942
$ws.Cells.Item(7, 12).Value = 0.257090072093213
$ws.Cells.Item(7, 14).Value = 2.286212572078823
$ws.Cells.Item(7, 15).Value = 5.174936621062329

$ws.Cells.Item(8, 2).Value = 0.904000713041512
$ws.Cells.Item(8, 3).Value = 0.110013124687498
$ws.Cells.Item(8, 4).Value = 0.1206250292390791
$ws.Cells.Item(8, 5).Value = 0.1322785252059617
$ws.Cells.Item(8, 6).Value = 1.942166452186164
$ws.Cells.Item(8, 10).Value = 0.1765307776212097
$ws.Cells.Item(8, 11).Value = 0.46654996979629
$ws.Cells.Item(8, 12).Value = 0.2618319617222582
$ws.Cells.Item(8, 14).Value = 2.259489058074048
$ws.Cells.Item(8, 15).Value = 5.140781823717333

$ws.Cells.Item(9, 2).Value = 0.9985964627636008
$ws.Cells.Item(9, 3).Value = 0.1125334978337804
$ws.Cells.Item(9, 4).Value = 0.1252209016905397
$ws.Cells.Item(9, 5).Value = 0.1323904215225511
$ws.Cells.Item(9, 6).Value = 1.930096956476284
$ws.Cells.Item(9, 10).Value = 0.1744829630837188
$ws.Cells.Item(9, 11).Value = 0.5520164383820827
$ws.Cells.Item(9, 12).Value = 0.2719334341492328
$ws.Cells.Item(9, 14).Value = 2.212136269059211
$ws.Cells.Item(9, 15).Value = 5.090868847179308

$ws.Cells.Item(10, 2).Value = 1.06974152430513
$ws.Cells.Item(10, 3).Value = 0.1143639177506444
$ws.Cells.Item(10, 4).Value = 0.128802893469711
$ws.Cells.Item(10, 5).Value = 0.1326821564595129
$ws.Cells.Item(10, 6).Value = 1.925113754240101
$ws.Cells.Item(10, 10).Value = 0.1732313271046309
$ws.Cells.Item(10, 11).Value = 0.6154363994871801
$ws.Cells.Item(10, 12).Value = 0.2798391457281326
$ws.Cells.Item(10, 14).Value = 2.180420048271353
$ws.Cells.Item(10, 15).Value = 5.064577286381677

$ws.Cells.Item(11, 2).Value = 1.102458693752141
$ws.Cells.Item(11, 3).Value = 0.1151918697627963
$ws.Cells.Item(11, 4).Value = 0.1304764475680713
$ws.Cells.Item(11, 5).Value = 0.1328600511805362
$ws.Cells.Item(11, 6).Value = 1.923687794362877
$ws.Cells.Item(11, 10).Value = 0.1727165747196153
$ws.Cells.Item(11, 11).Value = 0.6444200335887729
$ws.Cells.Item(11, 12).Value = 0.2835397588260236
$ws.Cells.Item(11, 14).Value = 2.16665771681649
$ws.Cells.Item(11, 15).Value = 5.054865695375895

$ws.Cells.Item(12, 2).Value = 1.114897916969227
$ws.Cells.Item(12, 3).Value = 0.1155047009852908
$ws.Cells.Item(12, 4).Value = 0.1311164599161998
$ws.Cells.Item(12, 5).Value = 0.1329338835726119
$ws.Cells.Item(12, 6).Value = 1.923268521099914
$ws.Cells.Item(12, 10).Value = 0.172529484977364
$ws.Cells.Item(12, 11).Value = 0.6554141001139726
$ws.Cells.Item(12, 12).Value = 0.2849559697561688
$ws.Cells.Item(12, 14).Value = 2.1615419278918
$ws.Cells.Item(12, 15).Value = 5.051511087788782

$ws.Cells.Item(13, 2).Value = 1.112216699417075
$ws.Cells.Item(13, 3).Value = 0.1154373584383777
$ws.Cells.Item(13, 4).Value = 0.1309783438316003
$ws.Cells.Item(13, 5).Value = 0.1329176952744682
$ws.Cells.Item(13, 6).Value = 1.923353454085117
$ws.Cells.Item(13, 10).Value = 0.1725694299011593
$ws.Cells.Item(13, 11).Value = 0.6530455140276672
$ws.Cells.Item(13, 12).Value = 0.2846503040646411
$ws.Cells.Item(13, 14).Value = 2.162639447867329
$ws.Cells.Item(13, 15).Value = 5.052219204819778

$ws.Cells.Item(14, 2).Value = 1.103481079871131
$ws.Cells.Item(14, 3).Value = 0.1152176206345956
$ws.Cells.Item(14, 4).Value = 0.1305289764124211
$ws.Cells.Item(14, 5).Value = 0.13286599596724
$ws.Cells.Item(14, 6).Value = 1.923650882597471
$ws.Cells.Item(14, 10).Value = 0.172701025800297
$ws.Cells.Item(14, 11).Value = 0.6453241529434592
$ws.Cells.Item(14, 12).Value = 0.2836559740643594
$ws.Cells.Item(14, 14).Value = 2.166234919632562
$ws.Cells.Item(14, 15).Value = 5.054583239345106

$ws.Cells.Item(15, 2).Value = 1.098136740470835
$ws.Cells.Item(15, 3).Value = 0.1150829336737047
$ws.Cells.Item(15, 4).Value = 0.1302545409571252
$ws.Cells.Item(15, 5).Value = 0.1328351700292068
$ws.Cells.Item(15, 6).Value = 1.923848779300158
$ws.Cells.Item(15, 10).Value = 0.172782651983276
$ws.Cells.Item(15, 11).Value = 0.640596999477907
$ws.Cells.Item(15, 12).Value = 0.2830488509347333
$ws.Cells.Item(15, 14).Value = 2.168449713990587
$ws.Cells.Item(15, 15).Value = 5.056073326548159

$ws.Cells.Item(16, 2).Value = 1.06761041293592
$ws.Cells.Item(16, 3).Value = 0.1143097126996864
$ws.Cells.Item(16, 4).Value = 0.1286944042884954
$ws.Cells.Item(16, 5).Value = 0.1326714372757891
$ws.Cells.Item(16, 6).Value = 1.925223847051029
$ws.Cells.Item(16, 10).Value = 0.1732660649201492
$ws.Cells.Item(16, 11).Value = 0.6135448861815007
$ws.Cells.Item(16, 12).Value = 0.2795993894301461
$ws.Cells.Item(16, 14).Value = 2.181332841567269
$ws.Cells.Item(16, 15).Value = 5.06525717613323

$ws.Cells.Item(17, 2).Value = 1.048973316645174
$ws.Cells.Item(17, 3).Value = 0.1138341460391388
$ws.Cells.Item(17, 4).Value = 0.1277485559650415
$ws.Cells.Item(17, 5).Value = 0.1325825438171009
$ws.Cells.Item(17, 6).Value = 1.926282632845499
$ws.Cells.Item(17, 10).Value = 0.173576599890648
$ws.Cells.Item(17, 11).Value = 0.5969830655724309
$ws.Cells.Item(17, 12).Value = 0.2775098730280234
$ws.Cells.Item(17, 14).Value = 2.189406684822848
$ws.Cells.Item(17, 15).Value = 5.071466812295171

$ws.Cells.Item(18, 2).Value = 1.038287019046692
$ws.Cells.Item(18, 3).Value = 0.1135601699427227
$ws.Cells.Item(18, 4).Value = 0.1272086827540591
$ws.Cells.Item(18, 5).Value = 0.1325356692679556
$ws.Cells.Item(18, 6).Value = 1.926970783114726
$ws.Cells.Item(18, 10).Value = 0.1737603539014714
$ws.Cells.Item(18, 11).Value = 0.5874697539451574
$ws.Cells.Item(18, 12).Value = 0.2763178598224982
$ws.Cells.Item(18, 14).Value = 2.194113188930113
$ws.Cells.Item(18, 15).Value = 5.075250099669944

$ws.Cells.Item(19, 2).Value = 1.034674560620715
$ws.Cells.Item(19, 3).Value = 0.1134673308596987
$ws.Cells.Item(19, 4).Value = 0.1270266063810084
$ws.Cells.Item(19, 5).Value = 0.1325205299923837
$ws.Cells.Item(19, 6).Value = 1.927217382201775
$ws.Cells.Item(19, 10).Value = 0.1738234536612531
$ws.Cells.Item(19, 11).Value = 0.5842508945162592
$ws.Cells.Item(19, 12).Value = 0.2759159552439741
$ws.Cells.Item(19, 14).Value = 2.19571748927908
$ws.Cells.Item(19, 15).Value = 5.076567421735319

$ws.Cells.Item(20, 2).Value = 1.050953830183715
$ws.Cells.Item(20, 3).Value = 0.1138848168754691
$ws.Cells.Item(20, 4).Value = 0.1278488137176055
$ws.Cells.Item(20, 5).Value = 0.1325915665564281
$ws.Cells.Item(20, 6).Value = 1.926161731614343
$ws.Cells.Item(20, 10).Value = 0.1735430108194933
$ws.Cells.Item(20, 11).Value = 0.5987447990605688
$ws.Cells.Item(20, 12).Value = 0.2777312900415438
$ws.Cells.Item(20, 14).Value = 2.188540727086463
$ws.Cells.Item(20, 15).Value = 5.070783881746451

$ws.Cells.Item(21, 2).Value = 1.106045594144035
$ws.Cells.Item(21, 3).Value = 0.1152821819945586
$ws.Cells.Item(21, 4).Value = 0.1306607967761693
$ws.Cells.Item(21, 5).Value = 0.132881006017211
$ws.Cells.Item(21, 6).Value = 1.923560246414411
$ws.Cells.Item(21, 10).Value = 0.172662160383485
$ws.Cells.Item(21, 11).Value = 0.6475916041455321
$ws.Cells.Item(21, 12).Value = 0.2839476303362432
$ws.Cells.Item(21, 14).Value = 2.165176244904252
$ws.Cells.Item(21, 15).Value = 5.053880103055064

$ws.Cells.Item(22, 2).Value = 1.14234196661414
$ws.Cells.Item(22, 3).Value = 0.1161913742754948
$ws.Cells.Item(22, 4).Value = 0.1325351277154567
$ws.Cells.Item(22, 5).Value = 0.1331078553316303
$ws.Cells.Item(22, 6).Value = 1.922563496479938
$ws.Cells.Item(22, 10).Value = 0.1721321381181475
$ws.Cells.Item(22, 11).Value = 0.6796239055656201
$ws.Cells.Item(22, 12).Value = 0.2880969787200627
$ws.Cells.Item(22, 14).Value = 2.150464064510033
$ws.Cells.Item(22, 15).Value = 5.044714795463818

$ws.Cells.Item(23, 2).Value = 1.122943575961983
$ws.Cells.Item(23, 3).Value = 0.1157064993367882
$ws.Cells.Item(23, 4).Value = 0.1315314406662083
$ws.Cells.Item(23, 5).Value = 0.1329833433204861
$ws.Cells.Item(23, 6).Value = 1.923031184679459
$ws.Cells.Item(23, 10).Value = 0.1724108489036347
$ws.Cells.Item(23, 11).Value = 0.6625179763571225
$ws.Cells.Item(23, 12).Value = 0.2858745095449962
$ws.Cells.Item(23, 14).Value = 2.158265185841299
$ws.Cells.Item(23, 15).Value = 5.049434388540476

$ws.Cells.Item(24, 2).Value = 1.050058351434927
$ws.Cells.Item(24, 3).Value = 0.113861910354295
$ws.Cells.Item(24, 4).Value = 0.1278034750116888
$ws.Cells.Item(24, 5).Value = 0.1325874741933823
$ws.Cells.Item(24, 6).Value = 1.926216143593479
$ws.Cells.Item(24, 10).Value = 0.1735581801672268
$ws.Cells.Item(24, 11).Value = 0.597948293451708
$ws.Cells.Item(24, 12).Value = 0.2776311585023024
$ws.Cells.Item(24, 14).Value = 2.188932024902792
$ws.Cells.Item(24, 15).Value = 5.071091970246869

$ws.Cells.Item(25, 2).Value = 0.9727146670326476
$ws.Cells.Item(25, 3).Value = 0.111855363086022
$ws.Cells.Item(25, 4).Value = 0.1239413156601188
$ws.Cells.Item(25, 5).Value = 0.1323232439732465
$ws.Cells.Item(25, 6).Value = 1.93267909891658
$ws.Cells.Item(25, 10).Value = 0.174992449605444
$ws.Cells.Item(25, 11).Value = 0.5287837424393729
$ws.Cells.Item(25, 12).Value = 0.269115322535697
$ws.Cells.Item(25, 14).Value = 2.224406232196885
$ws.Cells.Item(25, 15).Value = 5.102547192352432
